# Actualización automática del mapa (2025-10-07 09:02:34)
#
# The source "Caso 6243 / GARCIA, TEODORO 3252" record (row 11) has been
# resolved/removed, so the whole row is deleted and every following row
# shifts up by one. The "PD" (Q) and "N2" (R) columns are also dropped
# entirely from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete case row (row 11) - remaining rows shift up.
$ws.Rows(11).Delete()

# Drop the PD (Q) and N2 (R) columns entirely.
$ws.Columns("Q:R").Delete()
